$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.697207333333333
$ws.Range("H2").Value = 29.091622
$ws.Range("I2").Value = 0.3243108558382838
$ws.Range("J2").Value = 0.3243108558382838
$ws.Range("M2").Value = 670.6815796666667
$ws.Range("N2").Value = 2012.044739
$ws.Range("O2").Value = 0.8096423760738445
$ws.Range("P2").Value = 0.8096423760738443
$ws.Range("Q2").Value = 6503.738332675184
$ws.Range("R2").Value = 58533.64499407666
$ws.Range("S2").Value = 0.2625758119074502
$ws.Range("T2").Value = 0.2625758119074501
$ws.Range("G3").Value = 9.697207333333333
$ws.Range("H3").Value = 29.091622
$ws.Range("I3").Value = 0.3243108558382838
$ws.Range("J3").Value = 0.3243108558382838
$ws.Range("O3").Value = 0.09053284325954498
$ws.Range("P3").Value = 0.09053284325954496
$ws.Range("Q3").Value = 727.2370375775323
$ws.Range("R3").Value = 6545.133338197791
$ws.Range("S3").Value = 0.02936078387897624
$ws.Range("T3").Value = 0.02936078387897623
$ws.Range("G4").Value = 9.697207333333333
$ws.Range("H4").Value = 29.091622
$ws.Range("I4").Value = 0.3243108558382838
$ws.Range("J4").Value = 0.3243108558382838
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.114203333333333
$ws.Range("N4").Value = 3.34261
$ws.Range("O4").Value = 0.001345058909591271
$ws.Range("P4").Value = 0.00134505890959127
$ws.Range("Q4").Value = 10.80466073482445
$ws.Range("R4").Value = 97.24194661342001
$ws.Range("S4").Value = 0.0004362172061224539
$ws.Range("T4").Value = 0.0004362172061224537
$ws.Range("G5").Value = 9.697207333333333
$ws.Range("H5").Value = 29.091622
$ws.Range("I5").Value = 0.3243108558382838
$ws.Range("J5").Value = 0.3243108558382838
$ws.Range("M5").Value = 80.40286633333334
$ws.Range("N5").Value = 241.208599
$ws.Range("O5").Value = 0.09706180953056985
$ws.Range("P5").Value = 0.09706180953056984
$ws.Range("Q5").Value = 779.6832650286198
$ws.Range("R5").Value = 7017.149385257579
$ws.Range("S5").Value = 0.03147819851807161
$ws.Range("T5").Value = 0.03147819851807159
$ws.Range("G6").Value = 9.697207333333333
$ws.Range("H6").Value = 29.091622
$ws.Range("I6").Value = 0.3243108558382838
$ws.Range("J6").Value = 0.3243108558382838
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.174552666666667
$ws.Range("N6").Value = 3.523658
$ws.Range("O6").Value = 0.001417912226449558
$ws.Range("P6").Value = 0.001417912226449558
$ws.Range("Q6").Value = 11.38988073258622
$ws.Range("R6").Value = 102.508926593276
$ws.Range("S6").Value = 0.0004598443276634228
$ws.Range("T6").Value = 0.0004598443276634226
$ws.Range("I7").Value = 0.2826325233457075
$ws.Range("J7").Value = 0.2826325233457074
$ws.Range("M7").Value = 670.6815796666667
$ws.Range("N7").Value = 2012.044739
$ws.Range("O7").Value = 0.8096423760738445
$ws.Range("P7").Value = 0.8096423760738443
$ws.Range("Q7").Value = 5667.919969539306
$ws.Range("R7").Value = 51011.27972585375
$ws.Range("S7").Value = 0.228831267757365
$ws.Range("T7").Value = 0.2288312677573648
$ws.Range("I8").Value = 0.2826325233457075
$ws.Range("J8").Value = 0.2826325233457074
$ws.Range("O8").Value = 0.09053284325954498
$ws.Range("P8").Value = 0.09053284325954496
$ws.Range("S8").Value = 0.02558752593610663
$ws.Range("T8").Value = 0.02558752593610661
$ws.Range("I9").Value = 0.2826325233457075
$ws.Range("J9").Value = 0.2826325233457074
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.114203333333333
$ws.Range("N9").Value = 3.34261
$ws.Range("O9").Value = 0.001345058909591271
$ws.Range("P9").Value = 0.00134505890959127
$ws.Range("Q9").Value = 9.416115656950002
$ws.Range("R9").Value = 84.74504091255001
$ws.Range("S9").Value = 0.0003801573936664068
$ws.Range("T9").Value = 0.0003801573936664065
$ws.Range("I10").Value = 0.2826325233457075
$ws.Range("J10").Value = 0.2826325233457074
$ws.Range("M10").Value = 80.40286633333334
$ws.Range("N10").Value = 241.208599
$ws.Range("O10").Value = 0.09706180953056985
$ws.Range("P10").Value = 0.09706180953056984
$ws.Range("Q10").Value = 679.4834173400051
$ws.Range("R10").Value = 6115.350756060046
$ws.Range("S10").Value = 0.0274328241481254
$ws.Range("T10").Value = 0.02743282414812539
$ws.Range("I11").Value = 0.2826325233457075
$ws.Range("J11").Value = 0.2826325233457074
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.174552666666667
$ws.Range("N11").Value = 3.523658
$ws.Range("O11").Value = 0.001417912226449558
$ws.Range("P11").Value = 0.001417912226449558
$ws.Range("Q11").Value = 9.926126967710003
$ws.Range("R11").Value = 89.33514270939001
$ws.Range("S11").Value = 0.000400748110444169
$ws.Range("T11").Value = 0.0004007481104441687
$ws.Range("G12").Value = 3.910524
$ws.Range("H12").Value = 11.731572
$ws.Range("I12").Value = 0.1307825378608469
$ws.Range("J12").Value = 0.1307825378608469
$ws.Range("M12").Value = 670.6815796666667
$ws.Range("N12").Value = 2012.044739
$ws.Range("O12").Value = 0.8096423760738445
$ws.Range("P12").Value = 0.8096423760738443
$ws.Range("Q12").Value = 2622.716413644412
$ws.Range("R12").Value = 23604.44772279971
$ws.Range("S12").Value = 0.1058870847026236
$ws.Range("T12").Value = 0.1058870847026236
$ws.Range("G13").Value = 3.910524
$ws.Range("H13").Value = 11.731572
$ws.Range("I13").Value = 0.1307825378608469
$ws.Range("J13").Value = 0.1307825378608469
$ws.Range("O13").Value = 0.09053284325954498
$ws.Range("P13").Value = 0.09053284325954496
$ws.Range("Q13").Value = 293.26772042506
$ws.Range("R13").Value = 2639.40948382554
$ws.Range("S13").Value = 0.01184011500124156
$ws.Range("T13").Value = 0.01184011500124156
$ws.Range("G14").Value = 3.910524
$ws.Range("H14").Value = 11.731572
$ws.Range("I14").Value = 0.1307825378608469
$ws.Range("J14").Value = 0.1307825378608469
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 1.114203333333333
$ws.Range("N14").Value = 3.34261
$ws.Range("O14").Value = 0.001345058909591271
$ws.Range("P14").Value = 0.00134505890959127
$ws.Range("Q14").Value = 4.35711887588
$ws.Range("R14").Value = 39.21406988292
$ws.Range("S14").Value = 0.0001759102177686899
$ws.Range("T14").Value = 0.0001759102177686898
$ws.Range("G15").Value = 3.910524
$ws.Range("H15").Value = 11.731572
$ws.Range("I15").Value = 0.1307825378608469
$ws.Range("J15").Value = 0.1307825378608469
$ws.Range("M15").Value = 80.40286633333334
$ws.Range("N15").Value = 241.208599
$ws.Range("O15").Value = 0.09706180953056985
$ws.Range("P15").Value = 0.09706180953056984
$ws.Range("Q15").Value = 314.417338465292
$ws.Range("R15").Value = 2829.756046187628
$ws.Range("S15").Value = 0.01269398977977406
$ws.Range("T15").Value = 0.01269398977977406
$ws.Range("G16").Value = 3.910524
$ws.Range("H16").Value = 11.731572
$ws.Range("I16").Value = 0.1307825378608469
$ws.Range("J16").Value = 0.1307825378608469
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.174552666666667
$ws.Range("N16").Value = 3.523658
$ws.Range("O16").Value = 0.001417912226449558
$ws.Range("P16").Value = 0.001417912226449558
$ws.Range("Q16").Value = 4.593116392264001
$ws.Range("R16").Value = 41.338047530376
$ws.Range("S16").Value = 0.0001854381594389971
$ws.Range("T16").Value = 0.000185438159438997
$ws.Range("G17").Value = 3.066674
$ws.Range("H17").Value = 9.200022000000001
$ws.Range("I17").Value = 0.1025610400324547
$ws.Range("J17").Value = 0.1025610400324547
$ws.Range("M17").Value = 670.6815796666667
$ws.Range("N17").Value = 2012.044739
$ws.Range("O17").Value = 0.8096423760738445
$ws.Range("P17").Value = 0.8096423760738443
$ws.Range("Q17").Value = 2056.761762642696
$ws.Range("R17").Value = 18510.85586378426
$ws.Range("S17").Value = 0.08303776414448132
$ws.Range("T17").Value = 0.08303776414448127
$ws.Range("G18").Value = 3.066674
$ws.Range("H18").Value = 9.200022000000001
$ws.Range("I18").Value = 0.1025610400324547
$ws.Range("J18").Value = 0.1025610400324547
$ws.Range("O18").Value = 0.09053284325954498
$ws.Range("P18").Value = 0.09053284325954496
$ws.Range("Q18").Value = 229.9836270706434
$ws.Range("R18").Value = 2069.85264363579
$ws.Range("S18").Value = 0.00928514256179414
$ws.Range("T18").Value = 0.009285142561794135
$ws.Range("G19").Value = 3.066674
$ws.Range("H19").Value = 9.200022000000001
$ws.Range("I19").Value = 0.1025610400324547
$ws.Range("J19").Value = 0.1025610400324547
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 1.114203333333333
$ws.Range("N19").Value = 3.34261
$ws.Range("O19").Value = 0.001345058909591271
$ws.Range("P19").Value = 0.00134505890959127
$ws.Range("Q19").Value = 3.416898393046667
$ws.Range("R19").Value = 30.75208553742
$ws.Range("S19").Value = 0.0001379506406726002
$ws.Range("T19").Value = 0.0001379506406726001
$ws.Range("G20").Value = 3.066674
$ws.Range("H20").Value = 9.200022000000001
$ws.Range("I20").Value = 0.1025610400324547
$ws.Range("J20").Value = 0.1025610400324547
$ws.Range("M20").Value = 80.40286633333334
$ws.Range("N20").Value = 241.208599
$ws.Range("O20").Value = 0.09706180953056985
$ws.Range("P20").Value = 0.09706180953056984
$ws.Range("Q20").Value = 246.5693797099087
$ws.Range("R20").Value = 2219.124417389178
$ws.Range("S20").Value = 0.009954760132887269
$ws.Range("T20").Value = 0.009954760132887265
$ws.Range("G21").Value = 3.066674
$ws.Range("H21").Value = 9.200022000000001
$ws.Range("I21").Value = 0.1025610400324547
$ws.Range("J21").Value = 0.1025610400324547
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 1.174552666666667
$ws.Range("N21").Value = 3.523658
$ws.Range("O21").Value = 0.001417912226449558
$ws.Range("P21").Value = 0.001417912226449558
$ws.Range("Q21").Value = 3.601970124497334
$ws.Range("R21").Value = 32.41773112047601
$ws.Range("S21").Value = 0.0001454225526194001
$ws.Range("T21").Value = 0.0001454225526194001
$ws.Range("G22").Value = 4.775574
$ws.Range("H22").Value = 14.326722
$ws.Range("I22").Value = 0.1597130429227071
$ws.Range("J22").Value = 0.159713042922707
$ws.Range("M22").Value = 670.6815796666667
$ws.Range("N22").Value = 2012.044739
$ws.Range("O22").Value = 0.8096423760738445
$ws.Range("P22").Value = 0.8096423760738443
$ws.Range("Q22").Value = 3202.889514135062
$ws.Range("R22").Value = 28826.00562721556
$ws.Range("S22").Value = 0.1293104475619245
$ws.Range("T22").Value = 0.1293104475619244
$ws.Range("G23").Value = 4.775574
$ws.Range("H23").Value = 14.326722
$ws.Range("I23").Value = 0.1597130429227071
$ws.Range("J23").Value = 0.159713042922707
$ws.Range("O23").Value = 0.09053284325954498
$ws.Range("P23").Value = 0.09053284325954496
$ws.Range("Q23").Value = 358.14169679081
$ws.Range("R23").Value = 3223.27527111729
$ws.Range("S23").Value = 0.01445927588142642
$ws.Range("T23").Value = 0.01445927588142641
$ws.Range("G24").Value = 4.775574
$ws.Range("H24").Value = 14.326722
$ws.Range("I24").Value = 0.1597130429227071
$ws.Range("J24").Value = 0.159713042922707
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 1.114203333333333
$ws.Range("N24").Value = 3.34261
$ws.Range("O24").Value = 0.001345058909591271
$ws.Range("P24").Value = 0.00134505890959127
$ws.Range("Q24").Value = 5.32096046938
$ws.Range("R24").Value = 47.88864422442
$ws.Range("S24").Value = 0.0002148234513611202
$ws.Range("T24").Value = 0.0002148234513611201
$ws.Range("G25").Value = 4.775574
$ws.Range("H25").Value = 14.326722
$ws.Range("I25").Value = 0.1597130429227071
$ws.Range("J25").Value = 0.159713042922707
$ws.Range("M25").Value = 80.40286633333334
$ws.Range("N25").Value = 241.208599
$ws.Range("O25").Value = 0.09706180953056985
$ws.Range("P25").Value = 0.09706180953056984
$ws.Range("Q25").Value = 383.969837986942
$ws.Range("R25").Value = 3455.728541882479
$ws.Range("S25").Value = 0.01550203695171152
$ws.Range("T25").Value = 0.01550203695171152
$ws.Range("G26").Value = 4.775574
$ws.Range("H26").Value = 14.326722
$ws.Range("I26").Value = 0.1597130429227071
$ws.Range("J26").Value = 0.159713042922707
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 1.174552666666667
$ws.Range("N26").Value = 3.523658
$ws.Range("O26").Value = 0.001417912226449558
$ws.Range("P26").Value = 0.001417912226449558
$ws.Range("Q26").Value = 4.593116392264001
$ws.Range("R26").Value = 41.338047530376
$ws.Range("S26").Value = 0.0001854381594389971
$ws.Range("T26").Value = 0.000185438159438997
